$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a bare percentage ("NN%") must be force-typed
# as Text first -- otherwise Excel auto-parses the string into a percentage number
# (matching how typing "69%" into a General-formatted cell behaves in real Excel).
$percentTextCells = @("H3", "H12", "H17", "H24", "H25", "H27", "H31", "H36")
foreach ($ref in $percentTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$changes = [ordered]@{
    "E2" = "2026-02-06 17:47:44"
    "E3" = "2026-02-06 17:47:46"
    "H3" = "69%"
    "E4" = "2026-02-06 17:47:49"
    "J4" = "997.1 hPa"
    "E5" = "2026-02-06 17:47:51"
    "J5" = "997.4 hPa"
    "O5" = "11.2 °C"
    "E6" = "2026-02-06 17:47:53"
    "E7" = "2026-02-06 17:47:56"
    "J7" = "998.2 hPa"
    "E8" = "2026-02-06 17:47:58"
    "E9" = "2026-02-06 17:48:01"
    "O9" = "5.1 °C"
    "E10" = "2026-02-06 17:48:03"
    "O10" = "9.3 °C"
    "E11" = "2026-02-06 17:48:06"
    "E12" = "2026-02-06 17:48:08"
    "H12" = "57%"
    "E13" = "2026-02-06 17:48:11"
    "O13" = "10.4 °C"
    "E14" = "2026-02-06 17:48:13"
    "E15" = "2026-02-06 17:48:15"
    "J15" = "997.5 hPa"
    "O15" = "10.8 °C"
    "E16" = "2026-02-06 17:48:18"
    "E17" = "2026-02-06 17:48:20"
    "H17" = "83%"
    "J17" = "998.8 hPa"
    "O17" = "6.1 °C"
    "E18" = "2026-02-06 17:48:23"
    "N18" = "-6.0 °C 17:29 TU"
    "E19" = "2026-02-06 17:48:25"
    "I19" = "0.1 mm"
    "J19" = "999.7 hPa"
    "O19" = "10.0 °C"
    "E20" = "2026-02-06 17:48:28"
    "O20" = "-1.9 °C"
    "E21" = "2026-02-06 17:48:30"
    "J21" = "997.7 hPa"
    "O21" = "8.7 °C"
    "E22" = "2026-02-06 17:48:33"
    "E23" = "2026-02-06 17:48:35"
    "I23" = "0.2 mm"
    "J23" = "997.5 hPa"
    "O23" = "10.2 °C"
    "E24" = "2026-02-06 17:48:37"
    "H24" = "64%"
    "J24" = "996.9 hPa"
    "E25" = "2026-02-06 17:48:40"
    "H25" = "77%"
    "J25" = "998.3 hPa"
    "E26" = "2026-02-06 17:48:42"
    "E27" = "2026-02-06 17:48:45"
    "H27" = "79%"
    "J27" = "997.5 hPa"
    "O27" = "11.1 °C"
    "E28" = "2026-02-06 17:48:47"
    "J28" = "999.7 hPa"
    "O28" = "5.0 °C"
    "E29" = "2026-02-06 17:48:50"
    "E30" = "2026-02-06 17:48:52"
    "I30" = "0.1 mm"
    "E31" = "2026-02-06 17:48:54"
    "H31" = "82%"
    "J31" = "999.0 hPa"
    "E32" = "2026-02-06 17:48:56"
    "O32" = "15.7 °C"
    "E33" = "2026-02-06 17:48:59"
    "O33" = "10.4 °C"
    "E34" = "2026-02-06 17:49:01"
    "O34" = "8.9 °C"
    "E35" = "2026-02-06 17:49:04"
    "E36" = "2026-02-06 17:49:06"
    "H36" = "61%"
    "I36" = "1.2 mm"
    "J36" = "999.8 hPa"
    "O36" = "13.2 °C"
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}

